$d = $word.ActiveDocument

# Locate the paragraph that contains the astromap credit/link text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jenika*Hollana*CzechGlobe*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Select the paragraph's content, excluding the trailing paragraph mark.
    $r = $target.Range
    $r.End = $r.End - 1
    $r.Select() | Out-Null

    # Replace the whole (multi-run) paragraph content with a single run
    # containing the updated text (year changed from 2018 to 2022).
    $word.Selection.Delete() | Out-Null
    $word.Selection.Font.Reset() | Out-Null
    $word.Selection.TypeText(" Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).") | Out-Null
}
